$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (also updates the Table3 column names)
$ws.Range("J1").Value = 'human_chatbot_classification'
$ws.Range("M1").Value = 'human_complaint_classification'

# Fill in the human_complaint_classification column for rows 2-89
$ws.Range("M2").Value = 'blacklist'
$ws.Range("M3").Value = 'missing airtime'
$ws.Range("M4").Value = 'voice bundle'
$ws.Range("M5").Value = 'billing'
$ws.Range("M6").Value = 'data bundle'
$ws.Range("M7").Value = 'upgrade'
$ws.Range("M8").Value = 'fraud'
$ws.Range("M9").Value = 'missing airtime'
$ws.Range("M10").Value = 'network'
$ws.Range("M11").Value = 'internet'
$ws.Range("M12").Value = 'new contract'
$ws.Range("M13").Value = 'blacklist'
$ws.Range("M14").Value = 'billing'
$ws.Range("M15").Value = 'customer service'
$ws.Range("M16").Value = 'billing'
$ws.Range("M17").Value = 'subscription services'
$ws.Range("M18").Value = 'courier delivery'
$ws.Range("M19").Value = 'data bundle'
$ws.Range("M20").Value = 'customer service'
$ws.Range("M21").Value = 'customer service'
$ws.Range("M22").Value = 'missing airtime'
$ws.Range("M23").Value = 'number does not exist'
$ws.Range("M24").Value = 'billing'
$ws.Range("M25").Value = 'missing airtime'
$ws.Range("M26").Value = 'billing'
$ws.Range("M27").Value = 'customer service'
$ws.Range("M28").Value = 'internet'
$ws.Range("M29").Value = 'missing airtime'
$ws.Range("M30").Value = 'data bundle'
$ws.Range("M31").Value = 'cancel bundle'
$ws.Range("M32").Value = 'billing'
$ws.Range("M33").Value = 'downgrade account'
$ws.Range("M34").Value = 'customer service'
$ws.Range("M35").Value = 'new contract'
$ws.Range("M36").Value = 'upgrade'
$ws.Range("M37").Value = 'vodamail'
$ws.Range("M38").Value = 'courier delivery'
$ws.Range("M39").Value = 'billing'
$ws.Range("M40").Value = 'customer service'
$ws.Range("M41").Value = 'customer service'
$ws.Range("M42").Value = 'customer service'
$ws.Range("M43").Value = 'vodabucks'
$ws.Range("M44").Value = 'customer service'
$ws.Range("M45").Value = 'data bundle'
$ws.Range("M46").Value = 'data bundle'
$ws.Range("M47").Value = 'vodabucks'
$ws.Range("M48").Value = 'sim swap'
$ws.Range("M49").Value = 'data bundle'
$ws.Range("M50").Value = 'upgrade'
$ws.Range("M51").Value = 'sim swap'
$ws.Range("M52").Value = 'vodabucks'
$ws.Range("M53").Value = 'missing airtime'
$ws.Range("M54").Value = 'billing'
$ws.Range("M55").Value = 'customer service'
$ws.Range("M56").Value = 'courier delivery'
$ws.Range("M57").Value = 'sms spam'
$ws.Range("M58").Value = 'customer service'
$ws.Range("M59").Value = 'fraud'
$ws.Range("M60").Value = 'missing airtime'
$ws.Range("M61").Value = 'sms spam'
$ws.Range("M62").Value = 'customer service'
$ws.Range("M63").Value = 'billing'
$ws.Range("M64").Value = 'internet'
$ws.Range("M65").Value = 'fraud'
$ws.Range("M66").Value = 'subscription services'
$ws.Range("M67").Value = 'new contract'
$ws.Range("M68").Value = 'customer service'
$ws.Range("M69").Value = 'customer service'
$ws.Range("M70").Value = 'customer service'
$ws.Range("M71").Value = 'data bundle'
$ws.Range("M72").Value = 'billing'
$ws.Range("M73").Value = 'customer service'
$ws.Range("M74").Value = 'customer service'
$ws.Range("M75").Value = 'customer service'
$ws.Range("M76").Value = 'courier delivery'
$ws.Range("M77").Value = 'data bundle'
$ws.Range("M78").Value = 'fraud'
$ws.Range("M79").Value = 'rewards'
$ws.Range("M80").Value = 'customer service'
$ws.Range("M81").Value = 'sim swap'
$ws.Range("M82").Value = 'internet'
$ws.Range("M83").Value = 'customer service'
$ws.Range("M84").Value = 'billing'
$ws.Range("M85").Value = 'subscription services'
$ws.Range("M86").Value = 'data bundle'
$ws.Range("M87").Value = 'customer service'
$ws.Range("M88").Value = 'network'
$ws.Range("M89").Value = 'billing'

# Hide the now-superseded helper columns (chatbot_classification,
# human_chatbot_classification, chatbot_description) and widen the
# human_complaint_classification column so it is usable for review
$ws.Range("I1:K1").EntireColumn.Hidden = $true
$ws.Range("M1").EntireColumn.ColumnWidth = 19.166666666666668

# Row height touch-ups
$ws.Rows.Item(7).RowHeight = 409.5
$ws.Rows.Item(37).RowHeight = 245

# Move the view / selection to where the in-progress review is happening
$ws.Range("M90").Select()
